$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hras"
$ws.Range("C2").Value = "Agtr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.663483666666666
$ws.Range("H2").Value = 28.990451
$ws.Range("I2").Value = 0.4172798466714015
$ws.Range("J2").Value = 0.4172798466714016
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.212188333333333
$ws.Range("N2").Value = 6.636565
$ws.Range("O2").Value = 0.06691893508610254
$ws.Range("P2").Value = 0.06691893508610254
$ws.Range("Q2").Value = 21.37744582675722
$ws.Range("R2").Value = 192.397012440815
$ws.Range("S2").Value = 0.02792392297214234
$ws.Range("T2").Value = 0.02792392297214234

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hras"
$ws.Range("C3").Value = "Agtr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.663483666666666
$ws.Range("H3").Value = 28.990451
$ws.Range("I3").Value = 0.4172798466714015
$ws.Range("J3").Value = 0.4172798466714016
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.44817966666667
$ws.Range("N3").Value = 46.344539
$ws.Range("O3").Value = 0.4673090969404124
$ws.Range("P3").Value = 0.4673090969404123
$ws.Range("Q3").Value = 149.2832318885654
$ws.Range("R3").Value = 1343.549086997089
$ws.Range("S3").Value = 0.1949986683194464
$ws.Range("T3").Value = 0.1949986683194464

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hras"
$ws.Range("C4").Value = "Agtr1a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.663483666666666
$ws.Range("H4").Value = 28.990451
$ws.Range("I4").Value = 0.4172798466714015
$ws.Range("J4").Value = 0.4172798466714016
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.39736566666667
$ws.Range("N4").Value = 46.192097
$ws.Range("O4").Value = 0.4657719679734851
$ws.Range("P4").Value = 0.4657719679734851
$ws.Range("Q4").Value = 148.7921916295275
$ws.Range("R4").Value = 1339.129724665747
$ws.Range("S4").Value = 0.1943572553798128
$ws.Range("T4").Value = 0.1943572553798128

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Hras"
$ws.Range("C5").Value = "Agtr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.297564333333334
$ws.Range("H5").Value = 15.892693
$ws.Range("I5").Value = 0.2287546509102482
$ws.Range("J5").Value = 0.2287546509102482
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.212188333333333
$ws.Range("N5").Value = 6.636565
$ws.Range("O5").Value = 0.06691893508610254
$ws.Range("P5").Value = 0.06691893508610254
$ws.Range("Q5").Value = 11.71921001328278
$ws.Range("R5").Value = 105.472890119545
$ws.Range("S5").Value = 0.01530801763490695
$ws.Range("T5").Value = 0.01530801763490695

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hras"
$ws.Range("C6").Value = "Agtr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.297564333333334
$ws.Range("H6").Value = 15.892693
$ws.Range("I6").Value = 0.2287546509102482
$ws.Range("J6").Value = 0.2287546509102482
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.44817966666667
$ws.Range("N6").Value = 46.344539
$ws.Range("O6").Value = 0.4673090969404124
$ws.Range("P6").Value = 0.4673090969404123
$ws.Range("Q6").Value = 81.83772561705857
$ws.Range("R6").Value = 736.539530553527
$ws.Range("S6").Value = 0.1068991293377874
$ws.Range("T6").Value = 0.1068991293377874

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hras"
$ws.Range("C7").Value = "Agtr1a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.297564333333334
$ws.Range("H7").Value = 15.892693
$ws.Range("I7").Value = 0.2287546509102482
$ws.Range("J7").Value = 0.2287546509102482
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 15.39736566666667
$ws.Range("N7").Value = 46.192097
$ws.Range("O7").Value = 0.4657719679734851
$ws.Range("P7").Value = 0.4657719679734851
$ws.Range("Q7").Value = 81.56853518302457
$ws.Range("R7").Value = 734.1168166472211
$ws.Range("S7").Value = 0.1065475039375539
$ws.Range("T7").Value = 0.1065475039375539

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Hras"
$ws.Range("C8").Value = "Agtr1a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.197232333333334
$ws.Range("H8").Value = 24.591697
$ws.Range("I8").Value = 0.3539655024183503
$ws.Range("J8").Value = 0.3539655024183503
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 2.212188333333333
$ws.Range("N8").Value = 6.636565
$ws.Range("O8").Value = 0.06691893508610254
$ws.Range("P8").Value = 0.06691893508610254
$ws.Range("Q8").Value = 18.13382173342278
$ws.Range("R8").Value = 163.204395600805
$ws.Range("S8").Value = 0.02368699447905325
$ws.Range("T8").Value = 0.02368699447905325

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Hras"
$ws.Range("C9").Value = "Agtr1a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.197232333333334
$ws.Range("H9").Value = 24.591697
$ws.Range("I9").Value = 0.3539655024183503
$ws.Range("J9").Value = 0.3539655024183503
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 15.44817966666667
$ws.Range("N9").Value = 46.344539
$ws.Range("O9").Value = 0.4673090969404124
$ws.Range("P9").Value = 0.4673090969404123
$ws.Range("Q9").Value = 126.6323178547426
$ws.Range("R9").Value = 1139.690860692683
$ws.Range("S9").Value = 0.1654112992831786
$ws.Range("T9").Value = 0.1654112992831786

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Hras"
$ws.Range("C10").Value = "Agtr1a"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.197232333333334
$ws.Range("H10").Value = 24.591697
$ws.Range("I10").Value = 0.3539655024183503
$ws.Range("J10").Value = 0.3539655024183503
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 15.39736566666667
$ws.Range("N10").Value = 46.192097
$ws.Range("O10").Value = 0.4657719679734851
$ws.Range("P10").Value = 0.4657719679734851
$ws.Range("Q10").Value = 126.2157836909566
$ws.Range("R10").Value = 1135.942053218609
$ws.Range("S10").Value = 0.1648672086561184
$ws.Range("T10").Value = 0.1648672086561184
